# Applies hybrid bold + color (#2C3E50) highlighting to quantitative
# impact metrics (percentages, dollar amounts) across the resume body,
# splitting the affected runs exactly as Word would when you select the
# metric text and apply Bold + Font Color formatting to it.

$d = $word.ActiveDocument

# RGB(0x2C, 0x3E, 0x50) expressed the way Word's OM stores colors
# (R + G*256 + B*65536) -> renders as <w:color w:val="2C3E50"/>.
$metricColor = 5258796

function Set-MetricBold($Paragraph, $MetricText) {
    $rng = $Paragraph.Range.Duplicate
    $found = $rng.Find.Execute($MetricText, $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
    if ($found) {
        $rng.Font.Bold = $true
        $rng.Font.Color = $metricColor
    } else {
        Write-Output ("WARNING: metric not found -> " + $MetricText)
    }
    return $found
}

# --- Partner - Siege Analytics bullets ---------------------------------

# "...improving demographic classification accuracy from 23% to 64%"
$p = $d.Paragraphs(10)
Set-MetricBold $p '23%' | Out-Null
Set-MetricBold $p '64%' | Out-Null

# "...margin of error from ±4.2% to ±2.1%, increasing voter turnout
#  prediction accuracy from 71% to 87%, ..."
$p = $d.Paragraphs(12)
Set-MetricBold $p '±4.2%' | Out-Null
Set-MetricBold $p '±2.1%' | Out-Null
Set-MetricBold $p '71%' | Out-Null
Set-MetricBold $p '87%' | Out-Null

# "...reduced mapping costs by 73.5%, saving campaigns and organizations
#  $4.7M and enabling smaller nonprofits..."
$p = $d.Paragraphs(13)
Set-MetricBold $p '73.5%' | Out-Null
Set-MetricBold $p '$4.7M' | Out-Null

# "...valued over $2 trillion"
$p = $d.Paragraphs(14)
Set-MetricBold $p '$2' | Out-Null

# --- Data Products Manager - Helm/Murmuration bullet --------------------

# "...reducing processing time by 57%"
$p = $d.Paragraphs(19)
Set-MetricBold $p '57%' | Out-Null

# --- KEY ACHIEVEMENTS AND IMPACT bullets --------------------------------

# "...reducing mapping costs 73.5%"
$p = $d.Paragraphs(55)
Set-MetricBold $p '73.5%' | Out-Null

# "$4.7M savings enabled nonprofit access"
$p = $d.Paragraphs(56)
Set-MetricBold $p '$4.7M' | Out-Null

# "178% accuracy improvement in racial classification algorithms"
$p = $d.Paragraphs(58)
Set-MetricBold $p '178%' | Out-Null
